$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,13
$data[0,0] = 1.02
$data[0,1] = 1.038095141472792
$data[0,2] = 1.047694733142173
$data[0,3] = 1.036793429124184
$data[0,4] = 1.055374303016986
$data[0,5] = 1
$data[0,6] = $null
$data[0,7] = 1.035724590990446
$data[0,8] = 1.043194441642306
$data[0,9] = 1.0504566536221
$data[0,10] = 1.039586154858086
$data[0,11] = 1.058114936407015
$data[0,12] = 1.044675897405513
$data[1,0] = 1.02
$data[1,1] = 1.039119833424978
$data[1,2] = 1.048663467672557
$data[1,3] = 1.037666788216304
$data[1,4] = 1.056419720135004
$data[1,5] = 1
$data[1,6] = $null
$data[1,7] = 1.035878134986104
$data[1,8] = 1.043863407352827
$data[1,9] = 1.051236877508633
$data[1,10] = 1.040269019901001
$data[1,11] = 1.058973190397291
$data[1,12] = 1.045345813124074
$data[2,0] = 1.02
$data[2,1] = 1.039783168070326
$data[2,2] = 1.049290876425945
$data[2,3] = 1.038232517236614
$data[2,4] = 1.057096819229851
$data[2,5] = 1
$data[2,6] = $null
$data[2,7] = 1.035976093819804
$data[2,8] = 1.044295977195598
$data[2,9] = 1.051741688080043
$data[2,10] = 1.040710849206836
$data[2,11] = 1.059528577434564
$data[2,12] = 1.045778997265638
$data[3,0] = 1.02
$data[3,1] = 1.040062103162505
$data[3,2] = 1.049554775424602
$data[3,3] = 1.038470494237915
$data[3,4] = 1.057381625223256
$data[3,5] = 1
$data[3,6] = $null
$data[3,7] = 1.036016941524588
$data[3,8] = 1.044477758192573
$data[3,9] = 1.051953898685767
$data[3,10] = 1.040896586411224
$data[3,11] = 1.059762070887788
$data[3,12] = 1.045961036412479
$data[4,0] = 1.02
$data[4,1] = 1.040108941640916
$data[4,2] = 1.049599093201919
$data[4,3] = 1.038510460073247
$data[4,4] = 1.057429454385298
$data[4,5] = 1
$data[4,6] = $null
$data[4,7] = 1.036023780426309
$data[4,8] = 1.044508275820347
$data[4,9] = 1.051989529072233
$data[4,10] = 1.040927772026914
$data[4,11] = 1.059801275959265
$data[4,12] = 1.045991597378785
$data[5,0] = 1.02
$data[5,1] = 1.039786894943138
$data[5,2] = 1.04929440212243
$data[5,3] = 1.038235696530169
$data[5,4] = 1.057100624219086
$data[5,5] = 1
$data[5,6] = $null
$data[5,7] = 1.035976640942426
$data[5,8] = 1.04429840644331
$data[5,9] = 1.051744523695544
$data[5,10] = 1.04071333106805
$data[5,11] = 1.059531697351999
$data[5,12] = 1.045781429963161
$data[6,0] = 1.02
$data[6,1] = 1.038441380584394
$data[6,2] = 1.048022002469532
$data[6,3] = 1.037088458655606
$data[6,4] = 1.05572747278599
$data[6,5] = 1
$data[6,6] = $null
$data[6,7] = 1.035776770401347
$data[6,8] = 1.043420582619356
$data[6,9] = 1.050720343173228
$data[6,10] = 1.039816938128771
$data[6,11] = 1.058404978817931
$data[6,12] = 1.044902359528689
$data[7,0] = 1.02
$data[7,1] = 1.036072650739712
$data[7,2] = 1.045784291293456
$data[7,3] = 1.035071568024702
$data[7,4] = 1.05331277059079
$data[7,5] = 1
$data[7,6] = $null
$data[7,7] = 1.035413909678506
$data[7,8] = 1.041871505229953
$data[7,9] = 1.048915278360661
$data[7,10] = 1.038237177444352
$data[7,11] = 1.056419887259453
$data[7,12] = 1.043351082271711
$data[8,0] = 1.02
$data[8,1] = 1.034495018885999
$data[8,2] = 1.044295498132332
$data[8,3] = 1.03373017666971
$data[8,4] = 1.051706348322351
$data[8,5] = 1
$data[8,6] = $null
$data[8,7] = 1.035164854307805
$data[8,8] = 1.040837313481503
$data[8,9] = 1.047711717360484
$data[8,10] = 1.037183904706197
$data[8,11] = 1.055096752836951
$data[8,12] = 1.042315421852355
$data[9,0] = 1.02
$data[9,1] = 1.033812250385879
$data[9,2] = 1.043651557638314
$data[9,3] = 1.033150108978439
$data[9,4] = 1.051011559078505
$data[9,5] = 1
$data[9,6] = $null
$data[9,7] = 1.035055319378167
$data[9,8] = 1.040389153869927
$data[9,9] = 1.047190525456968
$data[9,10] = 1.036727809807282
$data[9,11] = 1.054523890058707
$data[9,12] = 1.041866625802726
$data[10,0] = 1.02
$data[10,1] = 1.033558693440418
$data[10,2] = 1.043412477764172
$data[10,3] = 1.032934761488393
$data[10,4] = 1.050753604493241
$data[10,5] = 1
$data[10,6] = $null
$data[10,7] = 1.035014379279644
$data[10,8] = 1.040222635633657
$data[10,9] = 1.046996925760081
$data[10,10] = 1.036558393252362
$data[10,11] = 1.054311113325651
$data[10,12] = 1.041699871091461
$data[11,0] = 1.02
$data[11,1] = 1.033613079835106
$data[11,2] = 1.043463756298281
$data[11,3] = 1.032980949018308
$data[11,4] = 1.050808931139593
$data[11,5] = 1
$data[11,6] = $null
$data[11,7] = 1.035023172556605
$data[11,8] = 1.040258356697604
$data[11,9] = 1.047038453807491
$data[11,10] = 1.036594733787697
$data[11,11] = 1.05435675419585
$data[11,12] = 1.041735642883416
$data[12,0] = 1.02
$data[12,1] = 1.03379129019402
$data[12,2] = 1.043631793008581
$data[12,3] = 1.033132305927663
$data[12,4] = 1.050990233998905
$data[12,5] = 1
$data[12,6] = $null
$data[12,7] = 1.035051940435577
$data[12,8] = 1.040375390474746
$data[12,9] = 1.04717452256142
$data[12,10] = 1.036713805832007
$data[12,11] = 1.054506301660369
$data[12,12] = 1.041852842861945
$data[13,0] = 1.02
$data[13,1] = 1.033901098593912
$data[13,2] = 1.043735340340364
$data[13,3] = 1.033225577239411
$data[13,4] = 1.051101956734553
$data[13,5] = 1
$data[13,6] = $null
$data[13,7] = 1.035069631633664
$data[13,8] = 1.040447491996473
$data[13,9] = 1.047258358255593
$data[13,10] = 1.036787169719506
$data[13,11] = 1.054598444130911
$data[13,12] = 1.041925046776103
$data[14,0] = 1.02
$data[14,1] = 1.034540339503552
$data[14,2] = 1.04433824949572
$data[14,3] = 1.033768690004797
$data[14,4] = 1.051752476144886
$data[14,5] = 1
$data[14,6] = $null
$data[14,7] = 1.035172088151051
$data[14,8] = 1.040867049069981
$data[14,9] = 1.047746306307306
$data[14,10] = 1.037214173828497
$data[14,11] = 1.055134773227938
$data[14,12] = 1.042345199668779
$data[15,0] = 1.02
$data[15,1] = 1.034941413963837
$data[15,2] = 1.044716630810265
$data[15,3] = 1.034109575267481
$data[15,4] = 1.052160744946571
$data[15,5] = 1
$data[15,6] = $null
$data[15,7] = 1.035235903550826
$data[15,8] = 1.041130133346177
$data[15,9] = 1.048052372323496
$data[15,10] = 1.037482017016711
$data[15,11] = 1.055471215577258
$data[15,12] = 1.042608657554824
$data[16,0] = 1.02
$data[16,1] = 1.035175388300337
$data[16,2] = 1.044937403346185
$data[16,3] = 1.034308481501379
$data[16,4] = 1.052398958594554
$data[16,5] = 1
$data[16,6] = $null
$data[16,7] = 1.035272962690113
$data[16,8] = 1.041283552369775
$data[16,9] = 1.048230891361914
$data[16,10] = 1.037638243369218
$data[16,11] = 1.055667462766284
$data[16,12] = 1.042762294451033
$data[17,0] = 1.02
$data[17,1] = 1.035255173332249
$data[17,2] = 1.045012692780527
$data[17,3] = 1.034376315852526
$data[17,4] = 1.052480196397859
$data[17,5] = 1
$data[17,6] = $null
$data[17,7] = 1.035285571190598
$data[17,8] = 1.041335858582231
$data[17,9] = 1.048291761031658
$data[17,10] = 1.037691512154711
$data[17,11] = 1.055734378998125
$data[17,12] = 1.042814674944311
$data[18,0] = 1.02
$data[18,1] = 1.034898378911987
$data[18,2] = 1.044676026922618
$data[18,3] = 1.034072993883765
$data[18,4] = 1.052116933561833
$data[18,5] = 1
$data[18,6] = $null
$data[18,7] = 1.035229073651441
$data[18,8] = 1.041101910363588
$data[18,9] = 1.048019534755974
$data[18,10] = 1.037453280185362
$data[18,11] = 1.055435117866708
$data[18,12] = 1.042580394492363
$data[19,0] = 1.02
$data[19,1] = 1.033738810202435
$data[19,2] = 1.043582307391172
$data[19,3] = 1.033087731896522
$data[19,4] = 1.050936841478957
$data[19,5] = 1
$data[19,6] = $null
$data[19,7] = 1.035043476020035
$data[19,8] = 1.040340928370783
$data[19,9] = 1.047134453861395
$data[19,10] = 1.036678742150305
$data[19,11] = 1.054462263378937
$data[19,12] = 1.041818331817842
$data[20,0] = 1.02
$data[20,1] = 1.033010053943221
$data[20,2] = 1.042895269387165
$data[20,3] = 1.032468927073248
$data[20,4] = 1.050195571529572
$data[20,5] = 1
$data[20,6] = $null
$data[20,7] = 1.034925314304643
$data[20,8] = 1.039862169298037
$data[20,9] = 1.046577934962226
$data[20,10] = 1.036191744556296
$data[20,11] = 1.053850649321034
$data[20,12] = 1.041338892852298
$data[21,0] = 1.02
$data[21,1] = 1.033396351942797
$data[21,2] = 1.043259421607639
$data[21,3] = 1.032796903601347
$data[21,4] = 1.050588466143552
$data[21,5] = 1
$data[21,6] = $null
$data[21,7] = 1.034988093220035
$data[21,8] = 1.040115996780363
$data[21,9] = 1.046872959160331
$data[21,10] = 1.03644991241964
$data[21,11] = 1.054174871816986
$data[21,12] = 1.041593080798764
$data[22,0] = 1.02
$data[22,1] = 1.034917824493125
$data[22,2] = 1.044694373858094
$data[22,3] = 1.034089523208962
$data[22,4] = 1.052136729801679
$data[22,5] = 1
$data[22,6] = $null
$data[22,7] = 1.035232160293569
$data[22,8] = 1.041114663218317
$data[22,9] = 1.048034372651873
$data[22,10] = 1.037466265128805
$data[22,11] = 1.055451428849544
$data[22,12] = 1.042593165457608
$data[23,0] = 1.02
$data[23,1] = 1.036684757245246
$data[23,2] = 1.046362265118693
$data[23,3] = 1.035592421792085
$data[23,4] = 1.053936436337893
$data[23,5] = 1
$data[23,6] = $null
$data[23,7] = 1.035508979348876
$data[23,8] = 1.042272240948187
$data[23,9] = 1.049381966031859
$data[23,10] = 1.038645603577923
$data[23,11] = 1.056933037961632
$data[23,12] = 1.04375238708063
$ws.Range("B2:N25").Value = $data